$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (even_MAG-GUT77590.fa) entirely; all rows below shift up by one,
# matching the rest of the diff (which is just the consequence of that shift)
# and the new dimension A1:D20.
$ws.Rows.Item(13).Delete()
